$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.408.67'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.671.57'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").Value = '''221.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").Value = '''0.5332'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("D8").Value = '''0.2661'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").Value = '''0.06388'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '''20.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.64%  '
$ws.Range("D11").Value = '''0.07860'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '''4.533'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '1.679.20'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").Value = '1.901.11'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '''0.5614'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("D16").Value = '0.0₅8189'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").Value = '''66.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = '26.405.54'
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").Value = '''197.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.82%  '
$ws.Range("D22").Value = '''10.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.26%  '
$ws.Range("D23").Value = '''6.074'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.24%  '
$ws.Range("D24").Value = '''1.011'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '''146.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''0.1229'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("D27").Value = '''7.243'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").Value = '''16.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").Value = '''1.501'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("D30").Value = '''0.05910'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.55%  '
$ws.Range("D31").Value = '''1.285'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").Value = '''3.559'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("D34").Value = '''1.617'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("D35").Value = '''0.9689'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("D36").Value = '''2.839'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("D37").Value = '''2.434'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("D38").Value = '''0.5835'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("D39").Value = '''0.01614'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("D40").Value = '1.080.51'
$ws.Range("E40").Value = '  +4.05%  '
$ws.Range("D41").Value = '''5.904'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("D42").Value = '''0.8657'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").Value = '''103.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("D45").Value = '1.811.02'
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("D46").Value = '''58.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = '''0.4411'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").Value = '''7.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("E51").Value = '  +0.23%  '
